# Trade #28 closed at 2026-02-18 00:13:22 - unknown UNKNOWN +0.000%
#
# This script:
#  1. Updates the Summary sheet's aggregate statistics.
#  2. Updates the Strategy Status sheet's "momentum" row.
#  3. Marks the open "momentum" DOWN trade (trade #58) as CLOSED with its
#     exit details, on both the "All Trades" sheet and the "momentum" sheet.
#  4. Appends a brand-new OPEN "momentum" DOWN trade (trade #87) to both the
#     "All Trades" sheet and the "momentum" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.77
$wsSummary.Range("B4").Value = 0.87
$wsSummary.Range("B5").Value = 0.31
$wsSummary.Range("B6").Value = 56
$wsSummary.Range("B8").Value = 22
$wsSummary.Range("B9").Value = 55.36

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C11").Value = 99.93000000000001
$wsStatus.Range("D11").Value = 1
$wsStatus.Range("E11").Value = -0.07000000000000001
$wsStatus.Range("F11").Value = -0.07000000000000001

# ---------------------------------------------------------------------------
# 3. All Trades sheet - close existing trade #58 (row 59) + append new trade
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close the open trade recorded in row 59 (Trade # 58, momentum, DOWN)
$wsAll.Range("G59").Value = 0.6899999999999999
$wsAll.Range("H59").Value = "CLOSED"
$wsAll.Range("I59").Value = -9.2105
$wsAll.Range("J59").Value = -0.07000000000000001
$wsAll.Range("K59").Value = 99.93000000000001
$wsAll.Range("L59").Value = "early_exit"
$wsAll.Range("M59").Value = 0.13

# Append new open trade (Trade # 87) as row 88
$wsAll.Range("A88").Value = 87

# Dates that look like "YYYY-MM-DD" get auto-converted to Excel date serials
# by plain .Value assignment. The source workbook stores these as literal
# text, so round-trip the value through a text-formula + paste-as-values so
# it lands as a plain string without leaving a formula or a cell style behind.
$wsAll.Range("B88").Formula = "=""2026-02-18"""
$wsAll.Range("B88").Copy()
$wsAll.Range("B88").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$wsAll.Range("C88").Value = "00:13:17"
$wsAll.Range("D88").Value = "momentum"
$wsAll.Range("E88").Value = "DOWN"
$wsAll.Range("F88").Value = 0.76
$wsAll.Range("G88").Value = ""
$wsAll.Range("H88").Value = "OPEN"
$wsAll.Range("I88").Value = 0
$wsAll.Range("J88").Value = 0
$wsAll.Range("K88").Value = 100
$wsAll.Range("L88").Value = ""
$wsAll.Range("M88").Value = 0
$wsAll.Range("N88").Value = 0
$wsAll.Range("O88").Value = 0
$wsAll.Range("P88").Value = 0.9
$wsAll.Range("Q88").Value = "Downward momentum: -1.980% over 10 samples"

# ---------------------------------------------------------------------------
# 4. "momentum" strategy sheet - close existing trade #58 (row 2) + append
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

# Close the open trade recorded in row 2 (Trade # 58, momentum, DOWN)
$wsMomentum.Range("G2").Value = 0.6899999999999999
$wsMomentum.Range("H2").Value = "CLOSED"
$wsMomentum.Range("I2").Value = -9.2105
$wsMomentum.Range("J2").Value = -0.07000000000000001
$wsMomentum.Range("K2").Value = 99.93000000000001
$wsMomentum.Range("P2").Value = "early_exit"
$wsMomentum.Range("Q2").Value = 0.13

# Append new open trade (Trade # 87) as row 18
$wsMomentum.Range("A18").Value = 87

$wsMomentum.Range("B18").Formula = "=""2026-02-18"""
$wsMomentum.Range("B18").Copy()
$wsMomentum.Range("B18").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$wsMomentum.Range("C18").Value = "00:13:17"
$wsMomentum.Range("D18").Value = "momentum"
$wsMomentum.Range("E18").Value = "DOWN"
$wsMomentum.Range("F18").Value = 0.76
$wsMomentum.Range("G18").Value = ""
$wsMomentum.Range("H18").Value = "OPEN"
$wsMomentum.Range("I18").Value = 0
$wsMomentum.Range("J18").Value = 0
$wsMomentum.Range("K18").Value = 100
$wsMomentum.Range("L18").Value = 0
$wsMomentum.Range("M18").Value = 0
$wsMomentum.Range("N18").Value = 0.9
$wsMomentum.Range("O18").Value = "Downward momentum: -1.980% over 10 samples"
$wsMomentum.Range("P18").Value = ""
$wsMomentum.Range("Q18").Value = 0

Write-Host "Edit complete"
